$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2021_13")
$ws.Range("C4").Value = "'1.3906"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'1.352"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.430"
$ws.Range("E4").Style = "Normal"
$ws.Range("C5").Value = "'0.9549"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'0.881"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.035"
$ws.Range("E5").Style = "Normal"
$ws.Range("C6").Value = "'1.1517"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'1.108"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.198"
$ws.Range("E6").Style = "Normal"
$ws.Range("C7").Value = "'1.6773"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'1.618"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.739"
$ws.Range("E7").Style = "Normal"
$ws.Range("C8").Value = "'1.5238"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'1.399"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.659"
$ws.Range("E8").Style = "Normal"
$ws.Range("C10").Value = "'1.5195"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'1.130"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.043"
$ws.Range("E10").Style = "Normal"
$ws.Range("C11").Value = "'1.6303"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.872"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.048"
$ws.Range("E11").Style = "Normal"
$ws.Range("C13").Value = "'7.9660"
$ws.Range("C13").Style = "Normal"
$ws.Range("E13").Value = "'79.660"
$ws.Range("E13").Style = "Normal"
$ws.Range("C16").Value = "'1.5633"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'1.513"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.615"
$ws.Range("E16").Style = "Normal"
$ws.Range("C17").Value = "'1.2579"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'1.182"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.339"
$ws.Range("E17").Style = "Normal"
$ws.Range("C18").Value = "'1.8112"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'1.746"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.879"
$ws.Range("E18").Style = "Normal"
$ws.Range("C19").Value = "'1.6564"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'1.567"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.751"
$ws.Range("E19").Style = "Normal"
$ws.Range("C20").Value = "'1.3345"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'1.202"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.481"
$ws.Range("E20").Style = "Normal"
$ws.Range("C22").Value = "'1.3414"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.960"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.875"
$ws.Range("E22").Style = "Normal"
$ws.Range("C23").Value = "'0.6520"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.436"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.974"
$ws.Range("E23").Style = "Normal"
$ws.Range("C25").Value = "'3.1674"
$ws.Range("C25").Style = "Normal"
$ws.Range("E25").Value = "'31.674"
$ws.Range("E25").Style = "Normal"

$ws = $wb.Worksheets.Item("2021_24")
$ws.Range("C4").Value = "'1.0948"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'1.065"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.126"
$ws.Range("E4").Style = "Normal"
$ws.Range("C5").Value = "'0.9336"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'0.844"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.033"
$ws.Range("E5").Style = "Normal"
$ws.Range("C6").Value = "'1.0134"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'0.959"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.071"
$ws.Range("E6").Style = "Normal"
$ws.Range("C7").Value = "'1.1255"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'1.074"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.179"
$ws.Range("E7").Style = "Normal"
$ws.Range("C8").Value = "'1.2537"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'1.187"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.324"
$ws.Range("E8").Style = "Normal"
$ws.Range("C9").Value = "'1.2049"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'1.110"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.308"
$ws.Range("E9").Style = "Normal"
$ws.Range("C10").Value = "'1.1353"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'1.007"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.280"
$ws.Range("E10").Style = "Normal"
$ws.Range("C11").Value = "'0.8069"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.645"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.009"
$ws.Range("E11").Style = "Normal"
$ws.Range("C13").Value = "'0.9435"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.383"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.327"
$ws.Range("E13").Style = "Normal"
$ws.Range("C16").Value = "'1.2577"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'1.231"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.285"
$ws.Range("E16").Style = "Normal"
$ws.Range("C17").Value = "'1.0684"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.135"
$ws.Range("E17").Style = "Normal"
$ws.Range("C18").Value = "'1.1528"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'1.112"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.195"
$ws.Range("E18").Style = "Normal"
$ws.Range("C19").Value = "'1.3206"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'1.273"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.370"
$ws.Range("E19").Style = "Normal"
$ws.Range("C20").Value = "'1.4842"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'1.405"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.568"
$ws.Range("E20").Style = "Normal"
$ws.Range("C21").Value = "'1.2668"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'1.142"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.405"
$ws.Range("E21").Style = "Normal"
$ws.Range("C22").Value = "'1.3137"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'1.111"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.554"
$ws.Range("E22").Style = "Normal"
$ws.Range("C23").Value = "'1.1608"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.853"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.580"
$ws.Range("E23").Style = "Normal"
$ws.Range("C25").Value = "'1.0200"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'0.414"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.516"
$ws.Range("E25").Style = "Normal"

$ws = $wb.Worksheets.Item("2022_06")
$ws.Range("C4").Value = "'1.0464"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.094"
$ws.Range("E4").Style = "Normal"
$ws.Range("C7").Value = "'0.9128"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'0.848"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.982"
$ws.Range("E7").Style = "Normal"
$ws.Range("C10").Value = "'0.8771"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.739"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.041"
$ws.Range("E10").Style = "Normal"
$ws.Range("C11").Value = "'0.9461"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.710"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.261"
$ws.Range("E11").Style = "Normal"
$ws.Range("C12").Value = "'0.9657"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.628"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.484"
$ws.Range("E12").Style = "Normal"
$ws.Range("C16").Value = "'1.0382"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'1.015"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.062"
$ws.Range("E16").Style = "Normal"
$ws.Range("C19").Value = "'0.9674"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'0.930"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.006"
$ws.Range("E19").Style = "Normal"
$ws.Range("C22").Value = "'1.0548"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.944"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.179"
$ws.Range("E22").Style = "Normal"
$ws.Range("C23").Value = "'0.9969"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.825"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.205"
$ws.Range("E23").Style = "Normal"
$ws.Range("C24").Value = "'1.1870"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'0.881"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.600"
$ws.Range("E24").Style = "Normal"
$ws.Range("C40").Value = "'1.6525"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'1.620"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.686"
$ws.Range("E40").Style = "Normal"
$ws.Range("C43").Value = "'1.7291"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'1.669"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.791"
$ws.Range("E43").Style = "Normal"
$ws.Range("C46").Value = "'1.4940"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.320"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.691"
$ws.Range("E46").Style = "Normal"
$ws.Range("C47").Value = "'1.0562"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.859"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'1.299"
$ws.Range("E47").Style = "Normal"
$ws.Range("C48").Value = "'1.2848"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.884"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1.868"
$ws.Range("E48").Style = "Normal"

$ws = $wb.Worksheets.Item("2022_47")
$ws.Range("C4").Value = "'1.2823"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'1.244"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.322"
$ws.Range("E4").Style = "Normal"
$ws.Range("C7").Value = "'1.2202"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'1.161"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.283"
$ws.Range("E7").Style = "Normal"
$ws.Range("C8").Value = "'1.4155"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'1.322"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.515"
$ws.Range("E8").Style = "Normal"
$ws.Range("C10").Value = "'1.8479"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'1.342"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.544"
$ws.Range("E10").Style = "Normal"
$ws.Range("C11").Value = "'0.6339"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.446"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.901"
$ws.Range("E11").Style = "Normal"
$ws.Range("C12").Value = "'0.5545"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.324"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.950"
$ws.Range("E12").Style = "Normal"
$ws.Range("C13").Value = "'3.9484"
$ws.Range("C13").Style = "Normal"
$ws.Range("E13").Value = "'39.484"
$ws.Range("E13").Style = "Normal"
